$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "skills" sheet: fix the coding-languages typo and expand the softwares
#    list.
# ---------------------------------------------------------------------------
$skills = $wb.Worksheets.Item("skills")

# Update "Softwares" row first (in-place text change, keeps its shared-string
# slot) then the "Coding languages" row (old string becomes unused and gets
# dropped, new one appended) - this ordering reproduces the same
# shared-string layout the original authoring session produced.
$skills.Range("B3").Value = "QGIS, Mendeley, Zotero, Inkscape, PhotoShop, GitHub"
$skills.Range("B2").Value = "R, UNIX, JavaScript, HTML, CSS"

# ---------------------------------------------------------------------------
# 2) Add a new "references" sheet at the end of the workbook.
# ---------------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$refs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$refs.Name = "references"

# Headers
$refs.Range("A1").Value = "names"
$refs.Range("C1").Value = "email"
$refs.Range("B1").Value = "institution"

# Names column (row order)
$refs.Range("A2").Value = "Prof. Marina B. Chiappero"
$refs.Range("A3").Value = "Prof. Esteban Soibelzon"
$refs.Range("A4").Value = "Prof. Sebastian Poljak"

# Institution column
$refs.Range("B2").Value = "Institute of Animal Diversity and Ecology (IDEA; UNC, CONICET)"
$refs.Range("B4").Value = "Southern Center for Scientific Research (CADIC; UNTDF, CONICET)"
$refs.Range("B3").Value = "Vertebrate Paleontology Division, La Plata Museum (UNLP, CONICET)"

# Email / contact column (row order)
$refs.Range("C2").Value = "Contact: marina.chiappero@gmail.com"
$refs.Range("C3").Value = "Contact: esoibel@gmail.com"
$refs.Range("C4").Value = "Contact: sebapoljak@hotmail.com"

# Header formatting - bold, matching the other sheets' header style
$refs.Range("A1:C1").Font.Bold = $true
$refs.Range("A1:C1").Font.Size = 12

# Column widths similar to the other sheets
$refs.Columns.Item(1).ColumnWidth = 21.83
$refs.Columns.Item(2).ColumnWidth = 19.67

# ---------------------------------------------------------------------------
# 3) Update cursor/selection position on a couple of sheets & land on the new
#    "references" sheet (mirrors the final state recorded in the workbook).
# ---------------------------------------------------------------------------
$education = $wb.Worksheets.Item("education")
$education.Activate()
$education.Range("A2").Select()

$skills.Activate()
$skills.Range("E9").Select()

$refs.Activate()
$refs.Range("D11").Select()

Write-Host "references sheet added"
